$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2445.971
$ws.Range("I138").Value = 3400
$ws.Range("J138").Value = 2320.8525
$ws.Range("K138").Value = 10200
$ws.Range("L138").Value = 6962.5575
$ws.Range("M138").Value = -5060
$ws.Range("N138").Value = -17242.5575

# Hunk 1: ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1692.3334
$ws.Range("I2").Value = 1581.5834
$ws.Range("J2").Value = 2135.3333
$ws.Range("K2").Value = 1581.5834
$ws.Range("L2").Value = 2135.3333
$ws.Range("M2").Value = -1468.5834
$ws.Range("N2").Value = -2361.3333

# Hunk 2: ARM row 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6999
$ws.Range("I6").Value = 6999
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 6999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -6826

# Hunk 3: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1783.25
$ws.Range("I61").Value = 1549.9286
$ws.Range("J61").Value = 3416.5
$ws.Range("K61").Value = 1549.9286
$ws.Range("L61").Value = 3416.5
$ws.Range("M61").Value = -1337.9286
$ws.Range("N61").Value = -3840.5

# Hunk 4: ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1157.55
$ws.Range("I97").Value = 1295.0625
$ws.Range("J97").Value = 607.5
$ws.Range("K97").Value = 1295.0625
$ws.Range("L97").Value = 607.5
$ws.Range("M97").Value = -799.0625
$ws.Range("N97").Value = -1599.5

# Hunk 5: ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1692.3334
$ws.Range("I116").Value = 1581.5834
$ws.Range("J116").Value = 2135.3333
$ws.Range("K116").Value = 1581.5834
$ws.Range("L116").Value = 2135.3333
$ws.Range("M116").Value = 712.4166
$ws.Range("N116").Value = -6723.3333

# Hunk 6: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20739.518
$ws.Range("I132").Value = 2248.2778
$ws.Range("J132").Value = 57722
$ws.Range("K132").Value = 6744.8334
$ws.Range("L132").Value = 173166
$ws.Range("M132").Value = -4214.8334
$ws.Range("N132").Value = -178226

# Hunk 7: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1783.25
$ws.Range("I136").Value = 1549.9286
$ws.Range("J136").Value = 3416.5
$ws.Range("K136").Value = 4649.7858
$ws.Range("L136").Value = 10249.5
$ws.Range("M136").Value = -2099.7858
$ws.Range("N136").Value = -15349.5

# Hunk 8: BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1692.3334
$ws.Range("I3").Value = 1581.5834
$ws.Range("J3").Value = 2135.3333
$ws.Range("K3").Value = 1581.5834
$ws.Range("L3").Value = 2135.3333
$ws.Range("M3").Value = -1467.5834
$ws.Range("N3").Value = -2363.3333

# Hunk 9: BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1322.5238
$ws.Range("I99").Value = 968.7143
$ws.Range("J99").Value = 2030.1428
$ws.Range("K99").Value = 968.7143
$ws.Range("L99").Value = 2030.1428
$ws.Range("M99").Value = 529.2857
$ws.Range("N99").Value = -5026.1428

# Hunk 10: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4390.815
$ws.Range("I134").Value = 4880.522
$ws.Range("J134").Value = 1575
$ws.Range("K134").Value = 14641.566
$ws.Range("L134").Value = 4725
$ws.Range("M134").Value = -12106.566
$ws.Range("N134").Value = -9795

# Hunk 11: CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()

# Hunk 12: CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1460
$ws.Range("I16").Value = 1552
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1552
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1265
$ws.Range("N16").Value = -1574

# Hunk 13: CRP row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 7998.5
$ws.Range("I25").Value = 7998.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 7998.5
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -7824.5

# Hunk 14: CRP row 92
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0

# Hunk 15: CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1460
$ws.Range("I113").Value = 1552
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1552
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 618
$ws.Range("N113").Value = -5340

# Hunk 16: CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 747.4
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 747.4
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2242.2
$ws.Range("N131").Value = -12322.2

# Hunk 17: GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2739.8
$ws.Range("I5").Value = 924.75
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 924.75
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -812.75
$ws.Range("N5").Value = -10224

# Hunk 18: GSM row 20
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 200000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 200000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 200000
$ws.Range("N20").Value = -200490

# Hunk 19: GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 855.625
$ws.Range("I97").Value = 790.8333
$ws.Range("J97").Value = 1050
$ws.Range("K97").Value = 790.8333
$ws.Range("L97").Value = 1050
$ws.Range("M97").Value = -294.8333
$ws.Range("N97").Value = -2042

# Hunk 20: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -550

# Hunk 21: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24211.76
$ws.Range("I132").Value = 4462.9473
$ws.Range("J132").Value = 86749.664
$ws.Range("K132").Value = 13388.8419
$ws.Range("L132").Value = 260248.992
$ws.Range("M132").Value = -10858.8419
$ws.Range("N132").Value = -265308.992

# Hunk 22: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5695.4
$ws.Range("I7").Value = 3531.3845
$ws.Range("J7").Value = 9714.286
$ws.Range("K7").Value = 3531.3845
$ws.Range("L7").Value = 9714.286
$ws.Range("M7").Value = -3419.3845
$ws.Range("N7").Value = -9938.286

# Hunk 23: LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 76.333336
$ws.Range("I55").Value = 80
$ws.Range("J55").Value = 75.111115
$ws.Range("K55").Value = 80
$ws.Range("L55").Value = 75.111115
$ws.Range("M55").Value = 93
$ws.Range("N55").Value = -421.111115

# Hunk 24: LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6425.4165
$ws.Range("I61").Value = 3566.6667
$ws.Range("J61").Value = 15001.667
$ws.Range("K61").Value = 3566.6667
$ws.Range("L61").Value = 15001.667
$ws.Range("M61").Value = -3364.6667
$ws.Range("N61").Value = -15405.667

# Hunk 25: LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3255.125
$ws.Range("I68").Value = 3200
$ws.Range("J68").Value = 3273.5
$ws.Range("K68").Value = 3200
$ws.Range("L68").Value = 3273.5
$ws.Range("M68").Value = -2451
$ws.Range("N68").Value = -4771.5

# Hunk 26: LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3255.125
$ws.Range("I71").Value = 3200
$ws.Range("J71").Value = 3273.5
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 16367.5
$ws.Range("M71").Value = -12256
$ws.Range("N71").Value = -23855.5

# Hunk 27: LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6425.4165
$ws.Range("I113").Value = 3566.6667
$ws.Range("J113").Value = 15001.667
$ws.Range("K113").Value = 3566.6667
$ws.Range("L113").Value = 15001.667
$ws.Range("M113").Value = -1396.6667
$ws.Range("N113").Value = -19341.667

# Hunk 28: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5695.4
$ws.Range("I126").Value = 3531.3845
$ws.Range("J126").Value = 9714.286
$ws.Range("K126").Value = 10594.1535
$ws.Range("L126").Value = 29142.858
$ws.Range("M126").Value = -8124.1535
$ws.Range("N126").Value = -34082.858

# Hunk 29: LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 47815.145
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 47815.145
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 47815.145
$ws.Range("N140").Value = -58175.145

# Hunk 30: WVR row 18
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 52000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 52000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 52000
$ws.Range("N18").Value = -52346

# Hunk 31: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2043.375
$ws.Range("I122").Value = 1946.4286
$ws.Range("J122").Value = 2722
$ws.Range("K122").Value = 5839.2858
$ws.Range("L122").Value = 8166
$ws.Range("M122").Value = -3389.2858
$ws.Range("N122").Value = -13066
